$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($col in @("C", "E", "G")) {
    # Mirror column A (values + styles) into the new column.
    $ws.Range("A1:A26").Copy($ws.Range($col + "1"))
}

# Rows 18, 23, 24 have a second (merged) cell in column B that must be
# mirrored too, then the copies re-merged.
$pairs = @(
    @{ src = "B18"; dstCol = "D" },
    @{ src = "B23:B24"; dstCol = "F" }
)

foreach ($col in @("D", "F", "H")) {
    $ws.Range("B18").Copy($ws.Range($col + "18"))
    $ws.Range("B23:B24").Copy($ws.Range($col + "23"))
}

foreach ($col in @("C", "E", "G")) {
    $next = [char]([int][char]$col + 1)
    $ws.Range($col + "18:" + $next + "18").Merge()
    $ws.Range($col + "23:" + $next + "24").Merge()
}

# Row 18 (A:H) was re-centered, creating a fresh alignment style shared by
# all eight cells across the merged blocks.
$ws.Range("A18:H18").HorizontalAlignment = -4108

# The new C/E/G11 cells got an explicit (distinct) center-alignment too.
$ws.Range("C11").HorizontalAlignment = -4108
$ws.Range("E11").HorizontalAlignment = -4108
$ws.Range("G11").HorizontalAlignment = -4108

$ws.Range("H7").Select()
